# ExpBoard_BoM.xlsx update: replace the placeholder/unconfirmed supplier
# codes with the exact parts used, and add Make/model + Supplier + Code
# detail for every line so the BoM is fully specified. The old
# "Total per board" / "Quantity required" / "Total cost (ex VAT)" summary
# block (column F) is removed, and three new LED lines (RED/GRN/BLUE) are
# appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the whole Price/Total column - it's not part of the new BoM.
$ws.Columns.Item(6).Delete()

# Clear out the old summary rows (Total per board / Quantity required /
# Total cost) - they lived in column F (now gone) and B, but B21/B23 must
# go too since the whole block is removed.
$ws.Range("A19:E19").ClearContents()
$ws.Range("A21:E21").ClearContents()
$ws.Range("A23:E23").ClearContents()

# --- Row 4: R1 (red LED limiting resistor) ---
$ws.Range("B4").Value = "Red LED current limiting resistor"
$ws.Range("C4").Value = "150R 0.25W"
$ws.Range("D4").Value = "Farnell"
$ws.Range("E4").Value = "'9339175 "

# --- Row 5: R2 (green LED limiting resistor) ---
$ws.Range("B5").Value = "Green LED current limiting resistor"
$ws.Range("C5").Value = "150R 0.25W"
$ws.Range("D5").Value = "Farnell"
$ws.Range("E5").Value = "'9339175 "

# --- Row 6: R3 (blue LED limiting resistor) ---
$ws.Range("B6").Value = "Blue LED current limiting resistor"
$ws.Range("C6").Value = "100R 0.25W"
$ws.Range("D6").Value = "Farnell"
$ws.Range("E6").Value = "'9339043"

# --- Row 7: R4 (switch pullup) ---
$ws.Range("C7").Value = "10K 0.25W"
$ws.Range("D7").Value = "Farnell"
$ws.Range("E7").Value = "'9339060 "

# --- Row 8: R5 (LDR potential divider) ---
$ws.Range("C8").Value = "10K 0.25W"
$ws.Range("D8").Value = "Farnell"
$ws.Range("E8").Value = "'9339060 "

# --- Row 9: C1 (decoupling capacitor) ---
$ws.Range("C9").Value = "0.1UF, 50V, Y5V"
$ws.Range("D9").Value = "Farnell"
$ws.Range("E9").Value = "9411887"

# --- Row 10: R-LDR ---
$ws.Range("E10").Value = "RE04698"

# --- Row 11: PIEZO ---
$ws.Range("C11").Value = "5V, 1mA"
$ws.Range("E11").Value = "LS03781"

# --- Row 12: TRI-LED (no code yet, still sourced from eBay) ---
# (no change)

# --- Row 13: TMP36 ---
$ws.Range("C13").Value = "TMP36GT9Z"

# --- Row 14: POT ---
$ws.Range("C14").Value = "10K"
$ws.Range("D14").Value = "Farnell"
$ws.Range("E14").Value = "'9608230"

# --- Row 15: SWITCH ---
$ws.Range("C15").Value = "Tactile"
$ws.Range("E15").Value = "1960939"

# --- Row 16: HEADER ---
$ws.Range("C16").Value = "10way, 1 row, Socket"
$ws.Range("E16").Value = "CN14535"

# --- New rows 17-19: discrete 5mm LEDs ---
$ws.Range("A17").Value = "RED"
$ws.Range("B17").Value = "Red LED 5mm"
$ws.Range("C17").Value = "T1 3/4, 20mA, 2V"
$ws.Range("D17").Value = "Farnell"
$ws.Range("E17").Value = "'1461624"

$ws.Range("A18").Value = "GRN"
$ws.Range("B18").Value = "Green 5mm LED"
$ws.Range("C18").Value = "T1 3/4, 20mA, 2V"
$ws.Range("D18").Value = "Farnell"
$ws.Range("E18").Value = "'1461633"

$ws.Range("A19").Value = "BLUE"
$ws.Range("B19").Value = "Blue 5mm LED"
$ws.Range("C19").Value = "T1 3/4, 20mA, 3.2V"
$ws.Range("D19").Value = "Farnell"
$ws.Range("E19").Value = "'1855507 "

# Column C needs to be a bit wider now it holds part numbers / values.
$ws.Columns.Item(3).ColumnWidth = 19.5703125

$ws.Range("C16").Select()
